$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: update title (D9) and link (E9)
$ws.Range("D9").Value = "MBA AI/BigData 및 MSc AI/DS 과정"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/mba-ai-bigdata-msc-ai-ds/#utm_source=rss&utm_medium=rss&utm_campaign=mba-ai-bigdata-msc-ai-ds"

# Row 32: update title (D32) and link (E32)
$ws.Range("D32").Value = "ROUGE : text summarization metric"
$ws.Range("E32").Value = "https://dodonam.tistory.com/368"
